$wb = $excel.ActiveWorkbook

# --- Sheet1: InvalidLoginTest -----------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Replace row 4 (kim124 -> peter) and drop the remaining kim12x rows.
$ws1.Range("A4").Value = "peter"
$ws1.Range("B4").Value = "peter123"
$ws1.Range("C4").Value = "Invalid credential"
$ws1.Range("A5:C15").Clear()

# --- Sheet2: rename Sheet2 -> AddValidEmployeeTest and add data -------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "AddValidEmployeeTest"

$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "First Name"
$ws2.Range("D1").Value = "Middle Name "
$ws2.Range("E1").Value = "Last Name"
$ws2.Range("F1").Value = "Expected Employee Name"

$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "admin123"
$ws2.Range("C2").Value = "Jack"
$ws2.Range("D2").Value = "wi"
$ws2.Range("E2").Value = "wick"
$ws2.Range("F2").Value = "Jack wick"

# Best-fit the data columns (B:F) like the recorded column widths.
$ws2.Columns.Item(2).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(3).ColumnWidth = 8.833333333333334
$ws2.Columns.Item(4).ColumnWidth = 11.5
$ws2.Columns.Item(5).ColumnWidth = 8.5
$ws2.Columns.Item(6).ColumnWidth = 21.666666666666668

# Make AddValidEmployeeTest the active tab with a fresh selection.
$ws2.Activate()
$ws2.Range("F3").Select()
